$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2464
$ws.Range("F4").Value = 32
$ws.Range("F5").Value = 1713
$ws.Range("F6").Value = 105
$ws.Range("F7").Value = 319
$ws.Range("F8").Value = 617
$ws.Range("F9").Value = 3538
$ws.Range("F10").Value = 928
$ws.Range("F11").Value = 1167
$ws.Range("F15").Value = 17
$ws.Range("F16").Value = 1306
$ws.Range("F17").Value = 1789
$ws.Range("F19").Value = 455
$ws.Range("F20").Value = 1548
$ws.Range("F21").Value = 8
$ws.Range("F22").Value = 1076
$ws.Range("F23").Value = 2261
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = 4290
$ws.Range("F28").Value = 2
$ws.Range("F31").Value = 1216
$ws.Range("F32").Value = 1

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F20").Value = 20
$ws.Range("F23").Value = 130
$ws.Range("F36").Value = 430
$ws.Range("F40").Value = 20
$ws.Range("F43").Value = 84

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2540
$ws.Range("F10").Value = 2972
$ws.Range("F11").Value = 488
$ws.Range("F12").Value = 805
$ws.Range("F13").Value = 214
$ws.Range("F14").Value = 214

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2540
$ws.Range("F6").Value = 2464
$ws.Range("F8").Value = 32
$ws.Range("F9").Value = 2972
$ws.Range("F10").Value = 488
$ws.Range("F11").Value = 805
$ws.Range("F12").Value = 214
$ws.Range("F15").Value = 1713
$ws.Range("F16").Value = 319
$ws.Range("F17").Value = 617
$ws.Range("F18").Value = 928
$ws.Range("F21").Value = 17
$ws.Range("F30").Value = 20
$ws.Range("F31").Value = 1789
$ws.Range("F33").Value = 1548
$ws.Range("F34").Value = 130
$ws.Range("F35").Value = 130
$ws.Range("F37").Value = 1076
$ws.Range("F40").Value = 2262
$ws.Range("F52").Value = 1216
